$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 122 (shifts old rows 122-179 down to 124-181)
$ws.Rows.Item(122).Insert()
$ws.Rows.Item(122).Insert()

# Populate new row 122
$ws.Cells.Item(122,1).Value = 8
$ws.Cells.Item(122,2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(122,3).Value = 'Coquimbo'
$ws.Cells.Item(122,4).Value = 44460
$ws.Cells.Item(122,5).Value = 4
$ws.Cells.Item(122,6).Value = 100112032
$ws.Cells.Item(122,7).Value = 'Zapallo italiano'
$ws.Cells.Item(122,8).Value = 'Sin especificar'
$ws.Cells.Item(122,9).Value = 'Primera'
$ws.Cells.Item(122,10).Value = 400
$ws.Cells.Item(122,11).Value = 10000
$ws.Cells.Item(122,12).Value = 11000
$ws.Cells.Item(122,13).Value = 10500
$ws.Cells.Item(122,14).Value = '$/caja 50 unidades'
$ws.Cells.Item(122,15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(122,16).Value = 210
$ws.Cells.Item(122,17).Value = 50
$ws.Cells.Item(122,18).Value = 'Hortaliza'

# Populate new row 123
$ws.Cells.Item(123,1).Value = 8
$ws.Cells.Item(123,2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(123,3).Value = 'Coquimbo'
$ws.Cells.Item(123,4).Value = 44460
$ws.Cells.Item(123,5).Value = 4
$ws.Cells.Item(123,6).Value = 100112032
$ws.Cells.Item(123,7).Value = 'Zapallo italiano'
$ws.Cells.Item(123,8).Value = 'Sin especificar'
$ws.Cells.Item(123,9).Value = 'Primera'
$ws.Cells.Item(123,10).Value = 600
$ws.Cells.Item(123,11).Value = 13000
$ws.Cells.Item(123,12).Value = 14000
$ws.Cells.Item(123,13).Value = 13500
$ws.Cells.Item(123,14).Value = '$/caja 70 unidades'
$ws.Cells.Item(123,15).Value = 'Provincia de Limarí'
$ws.Cells.Item(123,16).Value = 193
$ws.Cells.Item(123,17).Value = 70
$ws.Cells.Item(123,18).Value = 'Hortaliza'
